$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like bare percentages ("NN%") need an explicit
# text NumberFormat first, otherwise Excel auto-converts them to a numeric
# percentage (e.g. 0.76) instead of keeping the literal display string.
$percentCells = @("H3","H4","H5","H11","H12","H14","H15","H20","H21","H22","H26","H28","H29","H30","H32","H34")
foreach ($ref in $percentCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("E2").Value = "2026-02-06 04:17:41"
$ws.Range("O2").Value = "-1.4 °C"
$ws.Range("E3").Value = "2026-02-06 04:17:43"
$ws.Range("H3").Value = "76%"
$ws.Range("N3").Value = "-4.2 °C 3:31 TU"
$ws.Range("O3").Value = "-2.5 °C"
$ws.Range("E4").Value = "2026-02-06 04:17:46"
$ws.Range("H4").Value = "57%"
$ws.Range("J4").Value = "992.6 hPa"
$ws.Range("N4").Value = "11.2 °C 3:46 TU"
$ws.Range("O4").Value = "13.2 °C"
$ws.Range("E5").Value = "2026-02-06 04:17:48"
$ws.Range("H5").Value = "74%"
$ws.Range("J5").Value = "993.1 hPa"
$ws.Range("N5").Value = "6.3 °C 3:57 TU"
$ws.Range("O5").Value = "8.5 °C"
$ws.Range("E6").Value = "2026-02-06 04:17:51"
$ws.Range("J6").Value = "994.2 hPa"
$ws.Range("N6").Value = "13.9 °C 3:48 TU"
$ws.Range("O6").Value = "14.5 °C"
$ws.Range("E7").Value = "2026-02-06 04:17:53"
$ws.Range("J7").Value = "994.0 hPa"
$ws.Range("N7").Value = "9.5 °C 3:51 TU"
$ws.Range("O7").Value = "10.1 °C"
$ws.Range("E8").Value = "2026-02-06 04:17:55"
$ws.Range("N8").Value = "4.6 °C 3:47 TU"
$ws.Range("O8").Value = "6.2 °C"
$ws.Range("E9").Value = "2026-02-06 04:17:58"
$ws.Range("N9").Value = "1.2 °C 3:56 TU"
$ws.Range("O9").Value = "2.5 °C"
$ws.Range("E10").Value = "2026-02-06 04:18:00"
$ws.Range("N10").Value = "3.9 °C 3:55 TU"
$ws.Range("O10").Value = "5.3 °C"
$ws.Range("E11").Value = "2026-02-06 04:18:02"
$ws.Range("H11").Value = "83%"
$ws.Range("J11").Value = "994.9 hPa"
$ws.Range("N11").Value = "4.1 °C 3:30 TU"
$ws.Range("E12").Value = "2026-02-06 04:18:05"
$ws.Range("H12").Value = "59%"
$ws.Range("N12").Value = "10.3 °C 3:40 TU"
$ws.Range("O12").Value = "12.8 °C"
$ws.Range("E13").Value = "2026-02-06 04:18:07"
$ws.Range("E14").Value = "2026-02-06 04:18:09"
$ws.Range("H14").Value = "72%"
$ws.Range("I14").Value = "0.3 mm"
$ws.Range("O14").Value = "-3.6 °C"
$ws.Range("E15").Value = "2026-02-06 04:18:12"
$ws.Range("H15").Value = "81%"
$ws.Range("J15").Value = "993.2 hPa"
$ws.Range("N15").Value = "4.4 °C 3:59 TU"
$ws.Range("O15").Value = "7.7 °C"
$ws.Range("E16").Value = "2026-02-06 04:18:14"
$ws.Range("N16").Value = "3.4 °C 3:59 TU"
$ws.Range("O16").Value = "4.3 °C"
$ws.Range("E17").Value = "2026-02-06 04:18:17"
$ws.Range("J17").Value = "996.5 hPa"
$ws.Range("N17").Value = "2.6 °C 3:58 TU"
$ws.Range("O17").Value = "3.3 °C"
$ws.Range("E18").Value = "2026-02-06 04:18:19"
$ws.Range("N18").Value = "-5.2 °C 3:46 TU"
$ws.Range("E19").Value = "2026-02-06 04:18:21"
$ws.Range("J19").Value = "996.9 hPa"
$ws.Range("O19").Value = "6.3 °C"
$ws.Range("E20").Value = "2026-02-06 04:18:24"
$ws.Range("H20").Value = "72%"
$ws.Range("N20").Value = "-4.1 °C 3:50 TU"
$ws.Range("O20").Value = "-1.9 °C"
$ws.Range("E21").Value = "2026-02-06 04:18:26"
$ws.Range("H21").Value = "84%"
$ws.Range("J21").Value = "994.0 hPa"
$ws.Range("N21").Value = "3.0 °C 3:49 TU"
$ws.Range("O21").Value = "5.3 °C"
$ws.Range("E22").Value = "2026-02-06 04:18:29"
$ws.Range("H22").Value = "79%"
$ws.Range("N22").Value = "5.0 °C 3:59 TU"
$ws.Range("O22").Value = "8.9 °C"
$ws.Range("E23").Value = "2026-02-06 04:18:31"
$ws.Range("J23").Value = "993.4 hPa"
$ws.Range("L23").Value = "18.0 km/h - 45º 3:49 TU"
$ws.Range("M23").Value = "8.0 °C 3:37 TU"
$ws.Range("E24").Value = "2026-02-06 04:18:33"
$ws.Range("J24").Value = "992.3 hPa"
$ws.Range("E25").Value = "2026-02-06 04:18:35"
$ws.Range("J25").Value = "995.4 hPa"
$ws.Range("L25").Value = "10.4 km/h - 308º 3:50 TU"
$ws.Range("N25").Value = "1.6 °C 3:43 TU"
$ws.Range("O25").Value = "2.3 °C"
$ws.Range("E26").Value = "2026-02-06 04:18:38"
$ws.Range("H26").Value = "83%"
$ws.Range("N26").Value = "-1.0 °C 3:59 TU"
$ws.Range("E27").Value = "2026-02-06 04:18:40"
$ws.Range("J27").Value = "993.0 hPa"
$ws.Range("O27").Value = "7.7 °C"
$ws.Range("E28").Value = "2026-02-06 04:18:43"
$ws.Range("H28").Value = "88%"
$ws.Range("J28").Value = "995.9 hPa"
$ws.Range("N28").Value = "0.7 °C 3:43 TU"
$ws.Range("O28").Value = "3.3 °C"
$ws.Range("E29").Value = "2026-02-06 04:18:45"
$ws.Range("H29").Value = "59%"
$ws.Range("N29").Value = "8.1 °C 3:47 TU"
$ws.Range("O29").Value = "12.6 °C"
$ws.Range("E30").Value = "2026-02-06 04:18:47"
$ws.Range("H30").Value = "78%"
$ws.Range("L30").Value = "28.8 km/h - 345º 3:44 TU"
$ws.Range("O30").Value = "-3.5 °C"
$ws.Range("E31").Value = "2026-02-06 04:18:50"
$ws.Range("J31").Value = "996.4 hPa"
$ws.Range("E32").Value = "2026-02-06 04:18:52"
$ws.Range("H32").Value = "49%"
$ws.Range("J32").Value = "994.6 hPa"
$ws.Range("N32").Value = "13.4 °C 3:58 TU"
$ws.Range("O32").Value = "15.3 °C"
$ws.Range("E33").Value = "2026-02-06 04:18:55"
$ws.Range("O33").Value = "7.0 °C"
$ws.Range("E34").Value = "2026-02-06 04:18:57"
$ws.Range("H34").Value = "74%"
$ws.Range("N34").Value = "5.2 °C 3:53 TU"
$ws.Range("O34").Value = "8.8 °C"
$ws.Range("E35").Value = "2026-02-06 04:18:59"
$ws.Range("N35").Value = "-3.4 °C 3:32 TU"
$ws.Range("E36").Value = "2026-02-06 04:19:02"
$ws.Range("J36").Value = "996.1 hPa"
$ws.Range("O36").Value = "12.0 °C"
